$wb = $excel.ActiveWorkbook

# Overview sheet: update status columns (zh-cn, de-de) and latest handoff date
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-28-19 22:28:29"

# zh-cn sheet: update status and latest handoff datetime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-19 22:28:26"

# de-de sheet: update status and latest handoff datetime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-19 22:28:29"
